$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the "Tool" names first (column A) so the new shared-string table
# gets rBioNet before MetExplore, matching how the source workbook was
# authored (tool names typed down column A before the descriptions).
$ws.Range("A21").Value = "rBioNet"
$ws.Range("A20").Value = "MetExplore"

# Now the "Purpose of tool" column (column B).
$ws.Range("B20").Value = "Collaborative manual curation, visualization"
$ws.Range("B21").Value = "Assemble of reconstructions"

# Year of associated publications.
$ws.Range("C20").Value = 2018
$ws.Range("C21").Value = 2011

# Freely-available / Currently maintained.
$ws.Range("D20").Value = "Yes"
$ws.Range("D21").Value = "Yes"
$ws.Range("E20").Value = "Yes"
$ws.Range("E21").Value = "Yes"

# Copy the formatting of the last existing data row down onto the two new
# rows so the borders/fonts/alignment match the rest of the table.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A19:E19").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the (slightly taller) row height used for the new rows.
$ws.Rows.Item(20).RowHeight = 23.1
$ws.Rows.Item(21).RowHeight = 23.1

# The longer "Purpose of tool" text means column B needs to be widened.
$ws.Columns.Item(2).ColumnWidth = 37.32

# Restore the selection/active cell to the newly-entered cell.
$ws.Range("D20").Select()
